# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the Profits sheets (H:N columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7000
$ws.Range("I51").Value = 7000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = -6516
$ws.Range("M51").ClearContents()

$ws.Range("H123").Value = 4228293
$ws.Range("J123").Value = 73951.60000000001
$ws.Range("L123").Value = 73951.60000000001
$ws.Range("N123").Value = -83751.60000000001

$ws.Range("H124").Value = 59998.4
$ws.Range("J124").Value = 59998.4
$ws.Range("L124").Value = 59998.4
$ws.Range("N124").Value = -69818.39999999999

$ws.Range("H135").Value = 1328.091
$ws.Range("I135").Value = 1322.7188
$ws.Range("K135").Value = 11904.4692
$ws.Range("M135").Value = -9369.469200000001

$ws.Range("H138").Value = 2740
$ws.Range("I138").Value = 2045.619
$ws.Range("J138").Value = 3168.8823
$ws.Range("K138").Value = 6136.857
$ws.Range("L138").Value = 9506.6469
$ws.Range("M138").Value = -996.857
$ws.Range("N138").Value = -19786.6469

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3195.95
$ws.Range("I2").Value = 2663.5
$ws.Range("K2").Value = 2663.5
$ws.Range("M2").Value = -2550.5

$ws.Range("H14").Value = 1957.9375
$ws.Range("J14").Value = 830.1
$ws.Range("L14").Value = 830.1
$ws.Range("N14").Value = -1180.1

$ws.Range("H16").Value = 372.83334
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H32").Value = 2147.2446
$ws.Range("I32").Value = 2134.2717
$ws.Range("K32").Value = 2134.2717
$ws.Range("M32").Value = -1847.2717

$ws.Range("H110").Value = 859.15
$ws.Range("I110").Value = 799.1053000000001
$ws.Range("K110").Value = 799.1053000000001
$ws.Range("M110").Value = 1245.8947

$ws.Range("H116").Value = 3195.95
$ws.Range("I116").Value = 2663.5
$ws.Range("K116").Value = 2663.5
$ws.Range("M116").Value = -369.5

$ws.Range("H132").Value = 3021.8
$ws.Range("I132").Value = 1858.8
$ws.Range("K132").Value = 5576.4
$ws.Range("M132").Value = -3046.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3195.95
$ws.Range("I3").Value = 2663.5
$ws.Range("K3").Value = 2663.5
$ws.Range("M3").Value = -2549.5

$ws.Range("H86").Value = 2292.0908
$ws.Range("I86").Value = 2312.75
$ws.Range("J86").Value = 2237
$ws.Range("K86").Value = 2312.75
$ws.Range("L86").Value = 2237
$ws.Range("M86").Value = -1189.75
$ws.Range("N86").Value = -4483

$ws.Range("H89").Value = 2292.0908
$ws.Range("I89").Value = 2312.75
$ws.Range("J89").Value = 2237
$ws.Range("K89").Value = 11563.75
$ws.Range("L89").Value = 11185
$ws.Range("M89").Value = -5947.75
$ws.Range("N89").Value = -22417

$ws.Range("H99").Value = 31853.555
$ws.Range("I99").Value = 39811.715
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 39811.715
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -38313.715
$ws.Range("N99").Value = -6996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 7665.6665
$ws.Range("I19").Value = 6966.6665
$ws.Range("J19").Value = 8364.666999999999
$ws.Range("K19").Value = 6966.6665
$ws.Range("L19").Value = 8364.666999999999
$ws.Range("M19").Value = -6796.6665
$ws.Range("N19").Value = -8704.666999999999

$ws.Range("H24").Value = 7665.6665
$ws.Range("I24").Value = 6966.6665
$ws.Range("J24").Value = 8364.666999999999
$ws.Range("K24").Value = 6966.6665
$ws.Range("L24").Value = 8364.666999999999
$ws.Range("M24").Value = -6796.6665
$ws.Range("N24").Value = -8704.666999999999

$ws.Range("H31").Value = 5071.8184
$ws.Range("I31").Value = 1954
$ws.Range("K31").Value = 1954
$ws.Range("M31").Value = -1659

$ws.Range("H34").Value = 5071.8184
$ws.Range("I34").Value = 1954
$ws.Range("K34").Value = 1954
$ws.Range("M34").Value = -1752

$ws.Range("H58").Value = 1944.775
$ws.Range("I58").Value = 1514.2069
$ws.Range("K58").Value = 1514.2069
$ws.Range("M58").Value = -1311.2069

$ws.Range("H99").Value = 7629.353
$ws.Range("I99").Value = 4168.0557
$ws.Range("J99").Value = 11523.3125
$ws.Range("K99").Value = 4168.0557
$ws.Range("L99").Value = 11523.3125
$ws.Range("M99").Value = -2670.0557
$ws.Range("N99").Value = -14519.3125

$ws.Range("H126").Value = 7629.353
$ws.Range("I126").Value = 4168.0557
$ws.Range("J126").Value = 11523.3125
$ws.Range("K126").Value = 12504.1671
$ws.Range("L126").Value = 34569.9375
$ws.Range("M126").Value = -10034.1671
$ws.Range("N126").Value = -39509.9375

$ws.Range("H134").Value = 1564.5
$ws.Range("I134").Value = 1449.1765
$ws.Range("J134").Value = 1844.5714
$ws.Range("K134").Value = 4347.529500000001
$ws.Range("L134").Value = 5533.7142
$ws.Range("M134").Value = -1812.529500000001
$ws.Range("N134").Value = -10603.7142

$ws.Range("H136").Value = 1944.775
$ws.Range("I136").Value = 1514.2069
$ws.Range("K136").Value = 4542.620699999999
$ws.Range("M136").Value = -1992.620699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 157774820
$ws.Range("I4").Value = 168839380
$ws.Range("J4").Value = 25000000
$ws.Range("K4").Value = 506518140
$ws.Range("L4").Value = 75000000
$ws.Range("M4").Value = -506518028
$ws.Range("N4").Value = -75000224

$ws.Range("H12").Value = 28.875
$ws.Range("J12").Value = 45.25
$ws.Range("L12").Value = 135.75
$ws.Range("N12").Value = -481.75

$ws.Range("H21").Value = 1283.3334
$ws.Range("J21").Value = 1750
$ws.Range("L21").Value = 5250
$ws.Range("N21").Value = -5596

$ws.Range("H33").Value = 177.71428
$ws.Range("I33").Value = 163.5
$ws.Range("J33").Value = 196.66667
$ws.Range("K33").Value = 981
$ws.Range("L33").Value = 1180.00002
$ws.Range("M33").Value = -698
$ws.Range("N33").Value = -1746.00002

$ws.Range("H109").Value = 4666.5
$ws.Range("I109").Value = 2747.75
$ws.Range("J109").Value = 5306.0835
$ws.Range("K109").Value = 8243.25
$ws.Range("L109").Value = 15918.2505
$ws.Range("M109").Value = -7203.25
$ws.Range("N109").Value = -17998.2505

$ws.Range("H137").Value = 4221.4707
$ws.Range("I137").Value = 1994.5
$ws.Range("J137").Value = 6201
$ws.Range("K137").Value = 5983.5
$ws.Range("L137").Value = 18603
$ws.Range("M137").Value = -883.5
$ws.Range("N137").Value = -28803

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6352.375
$ws.Range("I80").Value = 5971
$ws.Range("J80").Value = 6842.7144
$ws.Range("K80").Value = 5971
$ws.Range("L80").Value = 6842.7144
$ws.Range("M80").Value = -4973
$ws.Range("N80").Value = -8838.714400000001

$ws.Range("H83").Value = 6352.375
$ws.Range("I83").Value = 5971
$ws.Range("J83").Value = 6842.7144
$ws.Range("K83").Value = 29855
$ws.Range("L83").Value = 34213.572
$ws.Range("M83").Value = -24863
$ws.Range("N83").Value = -44197.572

$ws.Range("H132").Value = 585324.25
$ws.Range("I132").Value = 1105902.9
$ws.Range("J132").Value = 12687.8
$ws.Range("K132").Value = 3317708.7
$ws.Range("L132").Value = 38063.39999999999
$ws.Range("M132").Value = -3315178.7
$ws.Range("N132").Value = -43123.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 80543.08
$ws.Range("J22").Value = 3645.4443
$ws.Range("L22").Value = 3645.4443
$ws.Range("N22").Value = -4235.4443

$ws.Range("H27").Value = 80543.08
$ws.Range("J27").Value = 3645.4443
$ws.Range("L27").Value = 3645.4443
$ws.Range("N27").Value = -3859.4443

$ws.Range("H46").Value = 981.3333
$ws.Range("I46").Value = 1052.909
$ws.Range("J46").Value = 784.5
$ws.Range("K46").Value = 1052.909
$ws.Range("L46").Value = 784.5
$ws.Range("M46").Value = -864.9090000000001
$ws.Range("N46").Value = -1160.5

$ws.Range("H82").Value = 1801.0667
$ws.Range("I82").Value = 1531.5714
$ws.Range("J82").Value = 2036.875
$ws.Range("K82").Value = 1531.5714
$ws.Range("L82").Value = 2036.875
$ws.Range("M82").Value = -1170.5714
$ws.Range("N82").Value = -2758.875

$ws.Range("H85").Value = 1801.0667
$ws.Range("I85").Value = 1531.5714
$ws.Range("J85").Value = 2036.875
$ws.Range("K85").Value = 1531.5714
$ws.Range("L85").Value = 2036.875
$ws.Range("M85").Value = -283.5714
$ws.Range("N85").Value = -4532.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13111160
$ws.Range("I136").Value = 18103760
$ws.Range("K136").Value = 54311280
$ws.Range("M136").Value = -54308730
